$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest reported period (column D: "6 ماهه منتهی به 1399/06") -- this
# shifts every later period one column to the left (D<-E, E<-F, ... L<-M).
$ws.Columns("D").Delete()

# Copy the formatting (styles/number formats) of the now-last data column (L)
# into the new trailing column (M) before we populate it with the newest
# reported period ("12 ماهه منتهی به 1401/12").
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$ws.Range("M1:M28").ColumnWidth = 28.1666666666667

# New period header + publish date for column M
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-27 (2)"

# Updated publish-date label for the 9-month 1401 period (column I after the shift)
$ws.Range("I9").Value = "1402-02-27 (9)"

# New cumulative income-statement figures for the newest period (column M)
$ws.Range("M11").Value = 79282474
$ws.Range("M12").Value = -25473270
$ws.Range("M13").Value = 53809204
$ws.Range("M14").Value = -3765171
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 3466328
$ws.Range("M17").Value = 53510361
$ws.Range("M18").Value = -166619
$ws.Range("M19").Value = -6761765
$ws.Range("M20").Value = 46581977
$ws.Range("M21").Value = 0
$ws.Range("M22").Value = 46581977
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 46581977
$ws.Range("M25").Value = 13199
$ws.Range("M26").Value = 3529200
$ws.Range("M27").Value = 13199
